# Auto-generated edit script: updates crypto price/volume table to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.774.86'
$ws.Range('E2').Value = '  +3.24%  '
$ws.Range('D3').Value = '2.551.88'
$ws.Range('E3').Value = '  +1.78%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '322.00'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '108.36'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.527'
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.556'
$ws.Range('E9').Value = '  +2.37%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.48'
$ws.Range('E10').Value = '  +1.46%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.43'
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.125'
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.24'
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('D15').Value = '2.948.69'
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('D16').Value = '2.583.14'
$ws.Range('E16').Value = '  +2.71%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.860'
$ws.Range('E17').Value = '  +1.60%  '
$ws.Range('D18').Value = '49.565.78'
$ws.Range('E18').Value = '  +3.20%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.15'
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '2.99'
$ws.Range('E20').Value = '  +9.58%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.69'
$ws.Range('E21').Value = '  +0.91%  '
$ws.Range('D22').Value = '0.0₃0942'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '284.49'
$ws.Range('E23').Value = '  +2.79%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '72.04'
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.53'
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '26.44'
$ws.Range('E26').Value = '  +2.05%  '
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.23'
$ws.Range('E28').Value = '  -7.10%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.144'
$ws.Range('E29').Value = '  +2.31%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.84'
$ws.Range('E30').Value = '  -2.36%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '35.40'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '49.53'
$ws.Range('E32').Value = '  +0.67%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.57'
$ws.Range('E33').Value = '  +1.07%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.37'
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0785'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('E37').Value = '  +2.43%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.65'
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.98'
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '120.88'
$ws.Range('E41').Value = '  -1.83%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.21'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '22.07'
$ws.Range('E43').Value = '  +1.68%  '
$ws.Range('E44').Value = '  +1.06%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.27'
$ws.Range('E45').Value = '  +4.11%  '
$ws.Range('D46').Value = '2.014.72'
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.00'
$ws.Range('E47').Value = '  +7.24%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.14'
$ws.Range('E48').Value = '  +7.58%  '
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '5.31'
$ws.Range('E50').Value = '  +2.43%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '81.20'
$ws.Range('E51').Value = '  +1.56%  '
